# controle_estoque.xlsx - "corrigindo matematica do total"
#
# Row 18 ("eli" / coca / 50) gets fixed to be the "renan" / coxinha / 10
# purchase, and the old row 19 ("renan" / coxinha / 10) is pushed down to
# make room for a handful of "pedro" / coca purchases that were missing
# from the ledger, finally landing (with name/product corrected back to
# "eli" / coca) as row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 fresh rows at row 19, shifting the old row 19 down to row 24
#    (this naturally carries its values/number-formats/styles with it).
$ws.Rows("19:23").Insert()

# 2) Correct row 18 in place (eli/coca/50 -> renan/coxinha/10).
$ws.Cells.Item(18, 2).Value = "renan"
$ws.Cells.Item(18, 3).Value = "coxinha"
$ws.Cells.Item(18, 4).Value = 10
$ws.Cells.Item(18, 6).Value = -10
$ws.Cells.Item(18, 10).Value = "'352456"

# 3) Fill in the newly inserted rows 19-23 with the missing "pedro" entries.
$ws.Cells.Item(19, 1).Value = 45769
$ws.Cells.Item(19, 2).Value = "pedro"
$ws.Cells.Item(19, 3).Value = "coca"
$ws.Cells.Item(19, 4).Value = 20
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = -20
$ws.Cells.Item(19, 7).Value = "aluno"
$ws.Cells.Item(19, 8).Value = "quarto"
$ws.Cells.Item(19, 9).Value = 25
$ws.Cells.Item(19, 10).Value = "'12"

$ws.Cells.Item(20, 1).Value = 45769
$ws.Cells.Item(20, 2).Value = "pedro"
$ws.Cells.Item(20, 3).Value = "coca"
$ws.Cells.Item(20, 4).Value = 20
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = -20
$ws.Cells.Item(20, 7).Value = "aluno"
$ws.Cells.Item(20, 8).Value = "quarto"
$ws.Cells.Item(20, 9).Value = 25
$ws.Cells.Item(20, 10).Value = "'12"

$ws.Cells.Item(21, 1).Value = 45769
$ws.Cells.Item(21, 2).Value = "pedro"
$ws.Cells.Item(21, 3).Value = "coca"
$ws.Cells.Item(21, 4).Value = 10
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = -10
$ws.Cells.Item(21, 7).Value = "aluno"
$ws.Cells.Item(21, 8).Value = "terceiro"
$ws.Cells.Item(21, 9).Value = 68765
$ws.Cells.Item(21, 10).Value = "asdasd"

$ws.Cells.Item(22, 1).Value = 45769
$ws.Cells.Item(22, 2).Value = "pedro"
$ws.Cells.Item(22, 3).Value = "coca"
$ws.Cells.Item(22, 4).Value = 10
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = -10
$ws.Cells.Item(22, 7).Value = "aluno"
$ws.Cells.Item(22, 8).Value = "quarto"
$ws.Cells.Item(22, 9).Value = 25
$ws.Cells.Item(22, 10).Value = "'12"

$ws.Cells.Item(23, 1).Value = 45769
$ws.Cells.Item(23, 2).Value = "pedro"
$ws.Cells.Item(23, 3).Value = "coca"
$ws.Cells.Item(23, 4).Value = 10
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = -10
$ws.Cells.Item(23, 7).Value = "aluno"
$ws.Cells.Item(23, 8).Value = "quarto"
$ws.Cells.Item(23, 9).Value = 25
$ws.Cells.Item(23, 10).Value = "'12"

# 4) Row 24 is the old row 19 (renan/coxinha/352456) shifted down; just fix
#    up the name/product/observacao back to eli/coca/"nao tem".
$ws.Cells.Item(24, 2).Value = "eli"
$ws.Cells.Item(24, 3).Value = "coca"
$ws.Cells.Item(24, 10).Value = "nao tem"
